# Add own student number: replace the placeholder tab + ellipsis
# ("…") that follows the first student number with the actual
# second student number "500799546" (no tab before it).

$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()

# Build the search string explicitly as a string concatenation
# (casting each char to [string] first) so PowerShell doesn't treat
# "+" between two [char] values as numeric addition.
$searchStr = [string][char]9 + [string][char]8230

$d.Content.Find.Execute($searchStr, $false, $false, $false, $false, $false, $true, 1, $false, "500799546", 2)
